# Update simulation results for the 380 kV case (pl_mw.xlsx line-loading data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated: B, C, D, E, G, H, I, J, K (F, L, M, N, O stay 0 / unchanged)
$columns = @(2, 3, 4, 5, 7, 8, 9, 10, 11)

# Each inner array holds the new values for row N, in the same order as $columns
$newValues = @(
    @(0.513310598074753, 0.1821040325354488, 0.6710255415307813, 0.2733999735695392, 2.151177162458538, 1.772764138966949, 1.506326244629058, 0.1420504691007523, 0.6903773610042947),  # row 2
    @(0.4839792836316974, 0.1736526897943236, 0.6614747676968307, 0.2686691276932507, 2.128169348583327, 1.767279396900619, 1.498284307995192, 0.1388827454660131, 0.6531332251788058),  # row 3
    @(0.4662991398375596, 0.1685817930317057, 0.6559622450545817, 0.2659144135013136, 2.114981480709361, 1.764548829094934, 1.493905714407234, 0.1370204881278596, 0.63071094826131),  # row 4
    @(0.4591772559090828, 0.1665450109327935, 0.6538042235892192, 0.2648295247606285, 2.109842794329566, 1.76359612932103, 1.492261717009043, 0.1362823590002691, 0.6216857603971278),  # row 5
    @(0.4579996833167002, 0.1662085930425263, 0.6534512218271686, 0.2646516541385395, 2.109003722738194, 1.763447594635807, 1.491997198521837, 0.1361610454044566, 0.6201939033227006),  # row 6
    @(0.4662027557677675, 0.1685542043416888, 0.655932783464948, 0.2658996298135605, 2.114911226017966, 1.764535332941421, 1.493882975155849, 0.1370104494844, 0.630588777668379),  # row 7
    @(0.5031288333053112, 0.1791654074357325, 0.6676594000024068, 0.2717376034728076, 2.143048772135074, 1.770740663892042, 1.503437141129844, 0.1409410345903979, 0.677443033801012),  # row 8
    @(0.5781553920847955, 0.2009176845921843, 0.693450626592238, 0.2843796473649149, 2.205711279327716, 1.787974309069227, 1.526625399537394, 0.1493082142258899, 0.7728682705735537),  # row 9
    @(0.6348798136689879, 0.2174838861122623, 0.7141134473945954, 0.2944014820436394, 2.256366747953706, 1.8037407965125, 1.546401152933825, 0.1558623067847975, 0.845156727969794),  # row 10
    @(0.661035538036316, 0.2251494631540254, 0.7238880556008951, 0.2991213904999768, 2.280426046045534, 1.811591542406489, 1.555998080628541, 0.158933348867194, 0.8785206470994638),  # row 11
    @(0.6709906329520265, 0.2280709843335273, 0.7276435041570153, 0.3009319247445461, 2.289683609102354, 1.814662246438559, 1.559718980496939, 0.1601092229557679, 0.8912238680651967),  # row 12
    @(0.6688443804040105, 0.2274409475578523, 0.7268322971845862, 0.3005409605724196, 2.287683285026588, 1.813996562124288, 1.558913754201484, 0.1598554010918605, 0.888484930812723),  # row 13
    @(0.6618535380891331, 0.2253894423876943, 0.724195935491565, 0.29926987874574, 2.281184724354233, 1.811842209745947, 1.556302460644787, 0.1590298291128107, 0.8795643648450948),  # row 14
    @(0.6575780190453031, 0.2241352791236295, 0.7225881252271904, 0.2984943287067097, 2.277223312064052, 1.810535349660512, 1.554714276047079, 0.1585258290153462, 0.8741092525134206),  # row 15
    @(0.633177543738924, 0.2169855382018397, 0.7134822074074521, 0.2940962688032371, 2.254814921983439, 1.803241404431901, 1.545786091559648, 0.1556634140564057, 0.8429859728429676),  # row 16
    @(0.6182986283159835, 0.212632654806157, 0.7079921306904851, 0.2914394651765875, 2.241328819963741, 1.79894074631099, 1.540463078203274, 0.1539303937992003, 0.8240157180569838),  # row 17
    @(0.6097737574581004, 0.2101411766986985, 0.7048696735388376, 0.2899264850104899, 2.233667544377141, 1.796530970377034, 1.537457954289494, 0.1529420318686761, 0.8131496606726785),  # row 18
    @(0.6068930673288264, 0.2092996946256562, 0.7038185214462658, 0.2894168147505383, 2.231089960247743, 1.795726020773117, 1.536450169804411, 0.1526088347877135, 0.809478349672645),  # row 19
    @(0.6198790902340932, 0.2130947647730466, 0.7085729057033348, 0.2917207187375794, 2.242754541666955, 1.79939194896329, 1.541023868635307, 0.1541140043296281, 0.8260304626444963),  # row 20
    @(0.6639055480668219, 0.2259915098939871, 0.7249688320682139, 0.2996425959115143, 2.283089518264774, 1.812472339224655, 1.557067103630047, 0.1592719680263173, 0.8821826772617669),  # row 21
    @(0.6929736556596424, 0.2345295233379829, 0.7359994137095782, 0.3049553027736778, 2.310306940166413, 1.821591275377955, 1.56805810508699, 0.1627184311761027, 0.9192837889896168),  # row 22
    @(0.6774325501080227, 0.2299625927911109, 0.7300833397655424, 0.3021074085974078, 2.295701900501172, 1.816672082139007, 1.562145601319486, 0.1608720666119439, 0.899445382100339),  # row 23
    @(0.619164472394516, 0.2128858104387348, 0.7083102318573538, 0.2915935189942616, 2.242109686169357, 1.799187764857692, 1.540770163525139, 0.154030969175011, 0.825119471422596),  # row 24
    @(0.557577958098733, 0.19493111317135, 0.6861731447443447, 0.280831225620652, 2.187952504851751, 1.782768145108378, 1.519873113929272, 0.1469735946839421, 0.7466718670111732)  # row 25
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNumber = 2 + $i
    $rowData = $newValues[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Cells.Item($rowNumber, $columns[$j]).Value = $rowData[$j]
    }
}
